$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 1129.1428
$ws.Range("I19").Value = 1400
$ws.Range("J19").Value = 1020.8
$ws.Range("K19").Value = 1400
$ws.Range("L19").Value = 1020.8
$ws.Range("M19").Value = -1225
$ws.Range("N19").Value = -1370.8
$ws.Range("H132").Value = 3944.0435
$ws.Range("I132").Value = 3668.818
$ws.Range("K132").Value = 11006.454
$ws.Range("M132").Value = -8476.454000000002
$ws.Range("H138").Value = 186221.81
$ws.Range("J138").Value = 287952.34
$ws.Range("L138").Value = 863857.02
$ws.Range("N138").Value = -874137.02

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H13").Value = 2800
$ws.Range("J13").Value = 2800
$ws.Range("L13").Value = 2800
$ws.Range("N13").Value = -3088
$ws.Range("H61").Value = 3566.9
$ws.Range("I61").Value = 3409.8572
$ws.Range("J61").Value = 3933.3333
$ws.Range("K61").Value = 3409.8572
$ws.Range("L61").Value = 3933.3333
$ws.Range("M61").Value = -3197.8572
$ws.Range("N61").Value = -4357.3333
$ws.Range("H62").Value = 30124.5
$ws.Range("J62").Value = 30124.5
$ws.Range("L62").Value = 30124.5
$ws.Range("N62").Value = -31372.5
$ws.Range("H65").Value = 30124.5
$ws.Range("J65").Value = 30124.5
$ws.Range("L65").Value = 90373.5
$ws.Range("N65").Value = -96613.5
$ws.Range("H136").Value = 3566.9
$ws.Range("I136").Value = 3409.8572
$ws.Range("J136").Value = 3933.3333
$ws.Range("K136").Value = 10229.5716
$ws.Range("L136").Value = 11799.9999
$ws.Range("M136").Value = -7679.571599999999
$ws.Range("N136").Value = -16899.9999

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H9").Value = 20000
$ws.Range("J9").Value = 20000
$ws.Range("L9").Value = 20000
$ws.Range("N9").Value = -20336
$ws.Range("H11").Value = 1834.6666
$ws.Range("J11").Value = 1752
$ws.Range("L11").Value = 1752
$ws.Range("N11").Value = -2032
$ws.Range("H12").Value = 575
$ws.Range("J12").Value = 600
$ws.Range("L12").Value = 600
$ws.Range("N12").Value = -936

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H5").Value = 1482.3334
$ws.Range("I5").Value = 378
$ws.Range("J5").Value = 2271.1428
$ws.Range("K5").Value = 378
$ws.Range("L5").Value = 2271.1428
$ws.Range("M5").Value = -266
$ws.Range("N5").Value = -2495.1428
$ws.Range("H13").Value = 31000
$ws.Range("I13").Value = 0
$ws.Range("K13").Value = 0
$ws.Range("M13").ClearContents()
$ws.Range("H22").Value = 404.16666
$ws.Range("I22").Value = 206.25
$ws.Range("J22").Value = 800
$ws.Range("K22").Value = 206.25
$ws.Range("L22").Value = 800
$ws.Range("M22").Value = 143.75
$ws.Range("N22").Value = -1500
$ws.Range("H53").Value = 26071.285
$ws.Range("J53").Value = 26071.285
$ws.Range("L53").Value = 26071.285
$ws.Range("N53").Value = -27285.285
$ws.Range("H63").Value = 53333.332
$ws.Range("J63").Value = 53333.332
$ws.Range("L63").Value = 53333.332
$ws.Range("N63").Value = -54705.332
$ws.Range("H66").Value = 53333.332
$ws.Range("J66").Value = 53333.332
$ws.Range("L66").Value = 159999.996
$ws.Range("N66").Value = -166863.996

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H132").Value = 3553.5881
$ws.Range("I132").Value = 2510.7693
$ws.Range("K132").Value = 22596.9237
$ws.Range("M132").Value = -20066.9237
$ws.Range("H134").Value = 4456.5186
$ws.Range("I134").Value = 2543.4736
$ws.Range("J134").Value = 9000
$ws.Range("K134").Value = 7630.4208
$ws.Range("L134").Value = 27000
$ws.Range("M134").Value = -2560.4208
$ws.Range("N134").Value = -37140

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H19").Value = 32001.2
$ws.Range("I19").Value = 0
$ws.Range("K19").Value = 0
$ws.Range("M19").ClearContents()

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 5140.731
$ws.Range("I22").Value = 900
$ws.Range("J22").Value = 8250.6
$ws.Range("K22").Value = 900
$ws.Range("L22").Value = 8250.6
$ws.Range("M22").Value = -605
$ws.Range("N22").Value = -8840.6
$ws.Range("H27").Value = 5140.731
$ws.Range("I27").Value = 900
$ws.Range("J27").Value = 8250.6
$ws.Range("K27").Value = 900
$ws.Range("L27").Value = 8250.6
$ws.Range("M27").Value = -793
$ws.Range("N27").Value = -8464.6
$ws.Range("H46").Value = 1575
$ws.Range("J46").Value = 1000
$ws.Range("L46").Value = 1000
$ws.Range("N46").Value = -1376
$ws.Range("H61").Value = 4358.65
$ws.Range("I61").Value = 4308.4287
$ws.Range("K61").Value = 4308.4287
$ws.Range("M61").Value = -4106.4287
$ws.Range("H62").Value = 77777
$ws.Range("J62").Value = 77777
$ws.Range("L62").Value = 77777
$ws.Range("N62").Value = -79025
$ws.Range("H65").Value = 77777
$ws.Range("J65").Value = 77777
$ws.Range("L65").Value = 233331
$ws.Range("N65").Value = -239571
$ws.Range("H68").Value = 3412.262
$ws.Range("I68").Value = 2294.95
$ws.Range("J68").Value = 4428
$ws.Range("K68").Value = 2294.95
$ws.Range("L68").Value = 4428
$ws.Range("M68").Value = -1545.95
$ws.Range("N68").Value = -5926
$ws.Range("H71").Value = 3412.262
$ws.Range("I71").Value = 2294.95
$ws.Range("J71").Value = 4428
$ws.Range("K71").Value = 11474.75
$ws.Range("L71").Value = 22140
$ws.Range("M71").Value = -7730.75
$ws.Range("N71").Value = -29628
$ws.Range("H93").Value = 6145.75
$ws.Range("I93").Value = 7503.647
$ws.Range("J93").Value = 2848
$ws.Range("K93").Value = 7503.647
$ws.Range("L93").Value = 2848
$ws.Range("M93").Value = -6255.647
$ws.Range("N93").Value = -5344
$ws.Range("H113").Value = 4358.65
$ws.Range("I113").Value = 4308.4287
$ws.Range("K113").Value = 4308.4287
$ws.Range("M113").Value = -2138.4287

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H3").Value = 32921.6
$ws.Range("J3").Value = 39652
$ws.Range("L3").Value = 39652
$ws.Range("N3").Value = -39880
$ws.Range("H49").Value = 70062
$ws.Range("J49").Value = 70062
$ws.Range("L49").Value = 70062
$ws.Range("N49").Value = -70522
$ws.Range("H126").Value = 1338.2174
$ws.Range("I126").Value = 1002.93335
$ws.Range("J126").Value = 1966.875
$ws.Range("K126").Value = 3008.80005
$ws.Range("L126").Value = 5900.625
$ws.Range("M126").Value = -538.8000499999998
$ws.Range("N126").Value = -10840.625
